$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I and J ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold + border + center/top alignment)
# from H1 onto the two new header cells so they share the same cell style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2-8 for columns I (I0) and J (IF) ---
$data = @(
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(5, 7),
    @(7, 7),
    @(2, 7),
    @(8, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
